# Add the new `summary_malfunction_reporting` field row to the
# device_classification_fields sheet (becomes new row 6; everything
# below shifts down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 6 (pushes old rows 6-18 to 7-19),
# inheriting the column B/C/D styles already used by the surrounding rows.
$ws.Rows("6:6").Insert()

$ws.Range("B6").Value = "summary_malfunction_reporting"
$ws.Range("C6").Value = "string"
$ws.Range("D6").Value = "The Voluntary Malfunction Summary Reporting Program allows participating companies to submit certain medical device malfunction reports in summary form on a quarterly basis.  The program applies to eligible devices regulated by the Center for Devices and Radiological Health (CDRH) and Center for Biologics Evaluation and Research (CBER), including device-led combination products.Value is one of the following:`nEligible = 510(K)`nIneligible = PMA"

# Match the row height Excel auto-computed for the wrapped description text.
$ws.Rows("6:6").RowHeight = 102

# Widen column B slightly now that it holds the longer field name, and
# drop its "best fit" auto-sizing (matches the checked-in width change,
# ~27.8 characters).
$ws.Columns("B").ColumnWidth = 27

# Reflect the reviewer's on-screen state when they made the edit: zoomed
# in on the new row, scrolled near the top, with D6 selected.
$ws.Application.ActiveWindow.Zoom = 170
$ws.Range("D6").Select()
